$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "46÷6=" },
    @{ Row = 1;  Col = 2; Text = "11÷8=" },
    @{ Row = 1;  Col = 3; Text = "72÷2=" },
    @{ Row = 1;  Col = 4; Text = "90÷6=" },
    @{ Row = 1;  Col = 5; Text = "14÷6=" },

    @{ Row = 5;  Col = 1; Text = "12÷2=" },
    @{ Row = 5;  Col = 2; Text = "88÷6=" },
    @{ Row = 5;  Col = 3; Text = "39÷2=" },
    @{ Row = 5;  Col = 4; Text = "91÷4=" },
    @{ Row = 5;  Col = 5; Text = "94÷6=" },

    @{ Row = 9;  Col = 1; Text = "27÷4=" },
    @{ Row = 9;  Col = 2; Text = "60÷8=" },
    @{ Row = 9;  Col = 3; Text = "43÷4=" },
    @{ Row = 9;  Col = 4; Text = "83÷7=" },
    @{ Row = 9;  Col = 5; Text = "61÷9=" },

    @{ Row = 13; Col = 1; Text = "61÷7=" },
    @{ Row = 13; Col = 2; Text = "27÷2=" },
    @{ Row = 13; Col = 3; Text = "20÷6=" },
    @{ Row = 13; Col = 4; Text = "65÷4=" },
    @{ Row = 13; Col = 5; Text = "20÷9=" },

    @{ Row = 17; Col = 1; Text = "64÷2=" },
    @{ Row = 17; Col = 2; Text = "30÷7=" },
    @{ Row = 17; Col = 3; Text = "48÷4=" },
    @{ Row = 17; Col = 4; Text = "23÷7=" },
    @{ Row = 17; Col = 5; Text = "72÷7=" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}
